# The document contains 8 occurrences of an "<id>...</id>" marker, each
# split across three runs: one run holding the literal text "<id>", a
# middle run holding the bare identifier (e.g. "p064r_1"), and a third
# run holding the literal text "</id>". The edit merges each triple into
# a single run (keeping the formatting of the "<id>" run) whose text is
# the full "<id>IDENT</id>" string. Two of the identifiers are also
# corrected (missing "r") while they're merged: p064_5 -> p064r_5 and
# p064_8 -> p064r_8.

$d = $word.ActiveDocument

$ids = @(
    @{old = "p064r_1"; new = "p064r_1"},
    @{old = "p064r_2"; new = "p064r_2"},
    @{old = "p064r_3"; new = "p064r_3"},
    @{old = "p064r_4"; new = "p064r_4"},
    @{old = "p064_5";  new = "p064r_5"},
    @{old = "p064r_6"; new = "p064r_6"},
    @{old = "p064r_7"; new = "p064r_7"},
    @{old = "p064_8";  new = "p064r_8"}
)

foreach ($item in $ids) {
    $old = $item.old
    $new = $item.new

    $rng = $d.Content
    $found = $rng.Find.Execute("<id>$old</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "WARNING: could not find <id>$old</id>"
        continue
    }

    $matchStart = $rng.Start
    $matchEnd = $rng.End
    # "<id>" is always 4 characters; everything after it (the identifier
    # plus the closing "</id>") is removed first so that re-inserting it
    # right after the opening run merges it into that run, picking up
    # the opening run's character formatting (Courier New etc.).
    $splitPoint = $matchStart + 4

    $tailRange = $d.Range($splitPoint, $matchEnd)
    $tailRange.Text = ""

    $headRange = $d.Range($matchStart, $splitPoint)
    $headRange.InsertAfter("$new</id>")
}
